# Edit: rebuild the "Requisitos" bullet list in LOB1223.docx
# (reordered, with some entries removed and some new ones added).
$d = $word.ActiveDocument

# Locate the paragraph that holds the Requisitos list: the one whose text
# begins with the old first requirement "LOB1257 ...".
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("LOB1257")) {
        $targetPara = $cand
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate the Requisitos list paragraph (LOB1257 ...)."
}

$newLines = @(
    'LOB1268 -  Leitura, Escrita e Comunicação Científica  (Requisito fraco)',
    'LOB1270 -  Química Experimental Aplicada  (Requisito fraco)',
    'LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)',
    'LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)',
    'LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)',
    'LOB1006 -  Cálculo IV  (Requisito fraco)',
    'LOB1040 -  Laboratório de Eletricidade  (Requisito fraco)',
    'LOB1053 -  Física III  (Requisito fraco)',
    'LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)',
    'LOB1003 -  Cálculo I  (Requisito fraco)',
    'LOB1011 -  Eletricidade Aplicada  (Requisito fraco)',
    'LOB1012 -  Estatística  (Requisito fraco)',
    'LOB1024 -  Mecânica  (Requisito fraco)',
    'LOB1036 -  Geometria Analítica  (Requisito fraco)',
    'LOB1037 -  Álgebra Linear  (Requisito fraco)',
    'LOB1038 -  Física Experimental I  (Requisito fraco)',
    'LOB1039 -  Física Experimental III  (Requisito fraco)',
    'LOB1041 -  Física Experimental II  (Requisito fraco)',
    'LOB1052 -  Cálculo III  (Requisito fraco)',
    'LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)',
    'LOB1004 -  Cálculo II  (Requisito fraco)',
    'LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito fraco)',
    'LOB1018 -  Física I  (Requisito fraco)',
    'LOB1019 -  Física II  (Requisito fraco)',
    'LOQ4233 -  Gestão de Negócios  (Requisito fraco)'
)

# Clear the paragraph's content but keep the paragraph mark itself.
$r = $targetPara.Range
$clearEnd = $r.End - 1
if ($clearEnd -gt $r.Start) {
    $d.Range($r.Start, $clearEnd).Text = ""
}

# Rebuild the list: one requirement per line, lines separated by a manual
# line break (w:br), matching the original formatting of the paragraph.
$pos = $targetPara.Range.Start
foreach ($line in $newLines) {
    $insLine = $d.Range($pos, $pos)
    $insLine.InsertAfter($line)
    $pos = $pos + $line.Length

    $insBreak = $d.Range($pos, $pos)
    $insBreak.InsertAfter([char]11)
    $pos = $pos + 1
}
